# Adds CypherOutput_Message, StatOutput, and StatOutput_Message sheets,
# mirroring the "Message" sheet layout for a new StatOutput cypher query
# and the stats results table (number_of_files / number_of_sample /
# number_of_cases / number_of_study).

$wb = $excel.ActiveWorkbook

# ---- Shared constants (same values already used on the "Message" sheet) ----
$neo4jUrlLabel = "Neo4j_URL:"
$neo4jUrl      = "bolt://ncias-q2251-c.nci.nih.gov:7687"
$userLabel     = "User_name:"
$userName      = "neo4j"
$pwdLabel      = "PWD:"
$pwd2          = "icdcDBneo4j0"
$cypherLabel   = "Cypher:"
$outputLabel   = "Output:"
$outputPath    = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC07_Canine_Filter_Diagnosis-LymphStg3_Neo4jData.xlsx'

$origCypher = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN [''Lymphoma :: Stage 3''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`'
$statCypher = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN[''Lymphoma :: Stage 3'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'

# ---- 1) CypherOutput_Message sheet: duplicate of the existing "Message" sheet ----
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cypherMsg = $wb.Worksheets.Add($null, $afterSheet)
$cypherMsg.Name = "CypherOutput_Message"
$cypherMsg.Range("A1").Value = $neo4jUrlLabel
$cypherMsg.Range("A2").Value = $neo4jUrl
$cypherMsg.Range("A3").Value = $userLabel
$cypherMsg.Range("A4").Value = $userName
$cypherMsg.Range("A5").Value = $pwdLabel
$cypherMsg.Range("A6").Value = $pwd2
$cypherMsg.Range("A7").Value = $cypherLabel
$cypherMsg.Range("A8").Value = $origCypher
$cypherMsg.Range("A9").Value = $outputLabel
$cypherMsg.Range("A10").Value = $outputPath

# ---- 2) StatOutput sheet: counts returned by the new stats cypher query ----
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$statOutput = $wb.Worksheets.Add($null, $afterSheet)
$statOutput.Name = "StatOutput"
$statOutput.Range("A1").Value = "number_of_files"
$statOutput.Range("B1").Value = "number_of_sample"
$statOutput.Range("C1").Value = "number_of_cases"
$statOutput.Range("D1").Value = "number_of_study"
$statOutput.Range("A2").Value = "'15"
$statOutput.Range("B2").Value = "'32"
$statOutput.Range("C2").Value = "'15"
$statOutput.Range("D2").Value = "'1"

# ---- 3) StatOutput_Message sheet: connection info repeated, then again
#         with the StatOutput cypher query ----
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$statMsg = $wb.Worksheets.Add($null, $afterSheet)
$statMsg.Name = "StatOutput_Message"
$statMsg.Range("A1").Value = $neo4jUrlLabel
$statMsg.Range("A2").Value = $neo4jUrl
$statMsg.Range("A3").Value = $userLabel
$statMsg.Range("A4").Value = $userName
$statMsg.Range("A5").Value = $pwdLabel
$statMsg.Range("A6").Value = $pwd2
$statMsg.Range("A7").Value = $cypherLabel
$statMsg.Range("A8").Value = $origCypher
$statMsg.Range("A9").Value = $outputLabel
$statMsg.Range("A10").Value = $outputPath
$statMsg.Range("A11").Value = $neo4jUrlLabel
$statMsg.Range("A12").Value = $neo4jUrl
$statMsg.Range("A13").Value = $userLabel
$statMsg.Range("A14").Value = $userName
$statMsg.Range("A15").Value = $pwdLabel
$statMsg.Range("A16").Value = $pwd2
$statMsg.Range("A17").Value = $cypherLabel
$statMsg.Range("A18").Value = $statCypher
$statMsg.Range("A19").Value = $outputLabel
$statMsg.Range("A20").Value = $outputPath

Write-Host "Added sheets:" $wb.Worksheets.Count
